# Generate Report for Handback
# Update the timestamp strings on the "Overview", "zh-cn", and "de-de"
# sheets to reflect the newly generated handback report times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 05:22:47"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 05:22:42"
$wsZhCn.Range("K2").Value = "2016-09-06 05:23:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-06 05:23:41"
